# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# 1) The second worker (row 63: 1047457297 / KATHERINE TATIANA DONCEL MANRIQUE / periodo 1612)
#    is removed from the account-statement table.
# 2) The remaining worker's period rows (16-62) are re-sorted in ascending
#    period order (1705 .. 2103) instead of descending, and the "Valor Mora"
#    (F) / "Salario Basico" (G) figures are refreshed.
# 3) The summary cells (worker count, period count, total Valor Mora) are
#    updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Give the (soon to be last) data row 62 the closing/bottom-border
#        look that row 63 currently has, by copying its formatting across,
#        then delete row 63 (the second worker) outright - this shifts the
#        trailing "firma" rows (68/69) up to (67/68) and fixes up the
#        merged-cell ranges automatically, just like a manual row delete
#        in Excel would.
$ws.Range("B63:J63").Copy()
$ws.Range("B62:J62").PasteSpecial(-4122)
$ws.Rows("63").Delete()

# --- 2. Rebuild the period / valor-mora / salario-basico columns for the
#        47 remaining rows (16-62), now in ascending chronological order.
$periods = @(
    "1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012",
    "2101","2102","2103"
)

$firstRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i
    if ($i -lt 16) {
        $mora = 29509
    } elseif ($i -lt 46) {
        $mora = 31249
    } else {
        $mora = 23958
    }
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $mora
    $ws.Range("G$row").Value = 781242
}

# --- 3. Refresh the summary header figures.
$ws.Range("C13").Value = 1        # Cant. Trabajadores
$ws.Range("F13").Value = 47       # Cant. Periodos
$ws.Range("E11").Value = 1433572  # VALOR MORA (total)

# --- 4. Column D ("Nombre Trabajador") no longer needs to fit the long
#        second worker's name, so it shrinks back down.
$ws.Columns("D:D").ColumnWidth = 26.83
